$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2..49. Column A holds line names like "C0_WA1_biomek-R1"
# which need their "C0_" prefix changed to "C0.1_". Column H holds the
# humidity value, changing from 90 to 85.
for ($row = 2; $row -le 49; $row++) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $cellA.Text -replace '^C0_', 'C0.1_'

    $ws.Cells.Item($row, 8).Value = 85
}
